# WRESBAL.xlsx update: append two new weekly observations to the "Data"
# sheet and refresh the FRED series metadata on "SeriesInfo" to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Data" sheet - append rows 110 and 111 (new weekly observations)
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

$data.Range("A110").Value = 45231
$data.Range("B110").Value = 3267.216
$data.Range("A111").Value = 45238
$data.Range("B111").Value = 3328.908

# Carry the existing date-column formatting (style used by every other
# row in column A) down onto the two new rows.
$data.Range("A109").Copy()
$data.Range("A110:A111").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. "SeriesInfo" sheet - refresh metadata fields reported by FRED
# ---------------------------------------------------------------------
$info = $wb.Worksheets.Item("SeriesInfo")

# These look like dates, so a plain assignment would get auto-converted
# to a date serial number by Excel. Lead with an apostrophe to force
# plain text, matching the original inline-string (non-date) storage.
$info.Range("B3").Value = "'2023-11-15"
$info.Range("B4").Value = "'2023-11-15"
$info.Range("B7").Value = "'2023-11-08"
$info.Range("B14").Value = "2023-11-09 15:39:01-06"
$info.Range("B15").Value = 73
